$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("OXM")

# Row 9 - Cost of Revenue: updated values
$ws.Range("D9").Value = 459700
$ws.Range("E9").Value = 425600
$ws.Range("F9").Value = 396700
$ws.Range("G9").Value = 388400
$ws.Range("H9").Value = 354500
$ws.Range("I9").Value = 369500
$ws.Range("J9").Value = 329100

# Row 10 - Gross Profit: updated values
$ws.Range("D10").Value = 626500
$ws.Range("E10").Value = 597000
$ws.Range("F10").Value = 572500
$ws.Range("G10").Value = 531900
$ws.Range("H10").Value = 495400
$ws.Range("I10").Value = 486000
$ws.Range("J10").Value = 429800

# Row 83 - Depreciation: J column now "NA"
$ws.Range("J83").Value = "NA"

# Row 94: J column now "NA"
$ws.Range("J94").Value = "NA"

# Row 100: J column now "NA"
$ws.Range("J100").Value = "NA"

# Row 101: J column now "NA"
$ws.Range("J101").Value = "NA"

$wb.Save()
